$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Taxa sheet: update the Notes text for the Didemnum/Trididemnum record.
# ---------------------------------------------------------------------------
$taxa = $wb.Worksheets.Item("Taxa")
$taxa.Range("AH5").Value = 'Reported as ''Didemnun / Trididemnum'' in PMLS records, this report includes all whitish, encrusting compound tunicates that often form large sheets. Field identification of these organisms is not possible.'

# ---------------------------------------------------------------------------
# 2. Materials sheet: insert three new Darwin Core columns (suborder,
#    infraorder, superfamily) right after the existing "order" column
#    (which pushes "family" and everything to its right three columns to
#    the right), then populate the new header/value cells and a few other
#    value cells that correspond to review-comment fixes.
# ---------------------------------------------------------------------------
$materials = $wb.Worksheets.Item("Materials")

# Insert 3 blank columns at AR:AT (immediately after "order" in AQ, before
# the old "family" column which was AR and becomes AU).
$materials.Range("AR1:AT2").EntireColumn.Insert()

# New column headers (row 1)
$materials.Range("AR1").Value = 'suborder'
$materials.Range("AS1").Value = 'infraorder'
$materials.Range("AT1").Value = 'superfamily'

# New column template values (row 2)
$materials.Range("AR2").Value = '${suborder}'
$materials.Range("AS2").Value = '${infraorder}'
$materials.Range("AT2").Value = '${superfamily}'

# scientificName template changed to use the occurrence summary taxon name
$materials.Range("AG2").Value = '${summary.taxonName}'

# scientificNameAuthorship (now shifted to BB) gets a template value
$materials.Range("BB2").Value = '${summary.Author}'

# eventTime (now shifted to EA) gets a template value
$materials.Range("EA2").Value = '!Date:HH:mm:ss'
